# Adds season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit "Created functions to get season record" which backfills
# each team's season W-L-T record alongside the existing player stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 56

$winsCol   = 30  # AD
$lossesCol = 31  # AE
$tiesCol   = 32  # AF

# --- Header row -------------------------------------------------------
$ws.Cells.Item(1, $winsCol).Value   = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value   = "Ties"

# Match the bold/centered/bordered header style used by the rest of row 1
# (copy format from the existing "Unnamed: 28" header in AC1).
$ws.Cells.Item(1, 29).Copy()
$ws.Range($ws.Cells.Item(1, $winsCol), $ws.Cells.Item(1, $tiesCol)).PasteSpecial(-4122)

# --- Data rows ---------------------------------------------------------
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, $winsCol).Value   = 98
    $ws.Cells.Item($row, $lossesCol).Value = 64
    $ws.Cells.Item($row, $tiesCol).Value   = 0
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
